# Actividad Asíncrona #14 -> #17, "Estructuras de Repetición" -> "Archivos en C"
#
# Original paragraph has two runs:
#   [bold]     "Actividad Asíncrona #14:"
#   [not bold] " Estructuras de Repetición"
#
# Target paragraph has five runs:
#   [bold]     "Actividad Asíncrona #1"
#   [bold]     "7"
#   [bold]     ":"
#   [not bold] " "
#   [not bold] "Archivos en C"

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "Actividad Asíncrona #14: Estructuras de Repetición",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target heading text."
}

$start = $rng.Start
$end = $rng.End

# Remove the old run contents so the new runs take their place.
$old = $d.Range($start, $end)
$old.Text = ""

$xmlFragment = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Actividad Asíncrona #1</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>7</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>:</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Archivos en C</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertionPoint = $d.Range($start, $start)
$insertionPoint.InsertXML($xmlFragment)

Write-Host "Updated heading to: $($d.Range($start, $start + 39).Text)"
